# Nuuchahnulth: update data progress
# Adds the 12/16 and 12/17 "Actual" progress-log entries (row 19 and row 20)
# to the Data sheet and moves the active selection to K21 (the next empty
# input cell), matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Row 19 (12/16/2019) actual counts ---
$ws.Range("J19").Value = 629
$ws.Range("K19").Value = 183
$ws.Range("L19").Value = 183
$ws.Range("M19").Value = 183
$ws.Range("N19").Value = 183
$ws.Range("O19").Formula = "=SUM((`$J19-`$D`$1), (`$K19-`$D`$1), (`$L19-`$D`$1), (`$M19-`$D`$1), (`$N19-`$D`$1))/(`$D`$3*5)"
$ws.Range("P19").Formula = "=1-O19"

# --- Row 20 (12/17/2019) actual counts ---
$ws.Range("J20").Value = 652
$ws.Range("K20").Value = 189
$ws.Range("L20").Value = 183
$ws.Range("M20").Value = 183
$ws.Range("N20").Value = 183
$ws.Range("O20").Formula = "=SUM((`$J20-`$D`$1), (`$K20-`$D`$1), (`$L20-`$D`$1), (`$M20-`$D`$1), (`$N20-`$D`$1))/(`$D`$3*5)"
$ws.Range("P20").Formula = "=1-O20"

# --- Move the active cell/selection to K21, the next empty entry cell ---
$ws.Activate()
$ws.Range("K21").Select()
